# Ordered seeding according to foreign key req. Updated Excel data
#
# The StudentTerm sheet's "B1" header cell used a stray shared string
# "StudentId" even though the rest of the workbook (Student, DegreePlan,
# etc.) consistently uses "StudentID". Correct that header so it points at
# the canonical "StudentID" text. Excel/the workbook engine will drop the
# now-unreferenced "StudentId" entry from the shared string table and
# renumber the remaining shared strings automatically on save, which also
# refreshes the dependent formula results on the StudentTerm sheet.
$wb = $excel.ActiveWorkbook

$studentTerm = $wb.Worksheets.Item("StudentTerm")
$studentTerm.Range("B1").Value = "StudentID"

# Restore the recorded selections (active cell) on the sheets whose view
# state changed.
$slot = $wb.Worksheets.Item("Slot")
[void]$slot.Range("F36").Select()

# StudentTerm remains the active/selected sheet, with B2 as its active cell.
[void]$studentTerm.Activate()
[void]$studentTerm.Range("B2").Select()
